# 6 hours by turn fix
# Shift the afternoon schedule earlier by 50 minutes, move lunch ("Almoço") to
# 12:20 (previously 13:00), and extend the timetable with two extra periods
# (17:30 and 18:20) so the turn covers a full 6 hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table for rows 1-17, columns A-F (row 1 header unchanged).
$data = @{
    2  = @("7:00","-","-","-","-","-")
    3  = @("7:50","-","-","-","João Rodrigues-CAD","Pedro Francisco-MTRM")
    4  = @("8:40","-","-","-","João Rodrigues-CAD","Pedro Francisco-MTRM")
    5  = @("9:30","Intervalo","Intervalo","Intervalo","Intervalo","Intervalo")
    6  = @("9:50","-","Euclides-Mecanica material","-","-","-")
    7  = @("10:40","-","Euclides-Mecanica material","-","-","-")
    8  = @("11:30","-","-","-","-","-")
    9  = @("12:20","Almoço","Almoço","Almoço","Almoço","Almoço")
    10 = @("13:00","-","-","-","-","-")
    11 = @("13:50","-","-","-","-","-")
    12 = @("14:40","-","-","-","-","-")
    13 = @("15:30","Intervalo","Intervalo","Intervalo","Intervalo","Intervalo")
    14 = @("15:50","-","-","-","-","-")
    15 = @("16:40","-","-","-","-","-")
    16 = @("17:30","-","-","-","-","-")
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}

# New row 17: only the time label is filled in, the rest of the row is blank.
$ws.Cells.Item(17, 1).Value = "18:20"
$ws.Cells.Item(17, 2).Value = ""
$ws.Cells.Item(17, 3).Value = ""
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = ""
$ws.Cells.Item(17, 6).Value = ""
